$d = $word.ActiveDocument

$replacements = @(
    @("18×11=", "57×18="),
    @("38×74=", "98×39="),
    @("24×81=", "36×26="),
    @("18×63=", "71×91="),
    @("72×17=", "71×70="),
    @("16×27=", "37×43="),
    @("81×78=", "94×96="),
    @("56×88=", "61×13="),
    @("31×39=", "40×18="),
    @("20×12=", "55×13="),
    @("98×87=", "28×73="),
    @("22×69=", "96×47="),
    @("80×90=", "33×69="),
    @("53×21=", "54×77="),
    @("13×73=", "49×39="),
    @("67×42=", "59×13="),
    @("58×93=", "53×94="),
    @("88×31=", "19×70="),
    @("51×33=", "51×28="),
    @("73×26=", "94×60="),
    @("77×96=", "90×24="),
    @("95×59=", "97×15="),
    @("77×30=", "55×88="),
    @("52×14=", "70×88="),
    @("65×80=", "11×77=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}
